$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: column D/E values forced to remain text (avoid Excel auto-numeric conversion)
$updates = @(
    @('D2', '26.648.87'),
    @('E2', '  +1.20%  '),
    @('D3', '1.633.27'),
    @('E4', '  +0.00%  '),
    @('D5', '213.11'),
    @('D6', '0.495'),
    @('E6', '  +1.48%  '),
    @('E7', '  +0.02%  '),
    @('E8', '  +1.01%  '),
    @('D9', '0.0626'),
    @('E9', '  +1.69%  '),
    @('D10', '19.04'),
    @('E10', '  +2.75%  '),
    @('E11', '  +3.48%  '),
    @('D12', '1.860.26'),
    @('D13', '1.632.77'),
    @('E13', '  +1.43%  '),
    @('E14', '  +1.40%  '),
    @('D15', '0.527'),
    @('E15', '  +1.98%  '),
    @('D16', '26.641.88'),
    @('E16', '  +1.31%  '),
    @('D17', '63.18'),
    @('E17', '  +1.26%  '),
    @('D18', '0.0₃0741'),
    @('E18', '  +1.71%  '),
    @('D19', '210.18'),
    @('E19', '  +3.71%  '),
    @('E20', '  -0.01%  '),
    @('D21', '4.31'),
    @('E21', '  +0.69%  '),
    @('E22', '  +1.07%  '),
    @('E23', '  +2.79%  '),
    @('D24', '1.92'),
    @('E24', '  +1.93%  '),
    @('D25', '147.05'),
    @('E25', '  +2.45%  '),
    @('E26', '  -0.01%  '),
    @('E27', '  -0.59%  '),
    @('D28', '6.88'),
    @('E28', '  +4.64%  '),
    @('D29', '15.41'),
    @('E29', '  +0.95%  '),
    @('D30', '0.0523'),
    @('E30', '  +5.23%  '),
    @('E31', '  -0.07%  '),
    @('E32', '  +1.59%  '),
    @('E33', '  +0.17%  '),
    @('D34', '1.51'),
    @('E34', '  +0.64%  '),
    @('D35', '2.35'),
    @('E35', '  -0.97%  '),
    @('B36', 'VeChain'),
    @('C36', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D36', '0.0171'),
    @('E36', '  +2.08%  '),
    @('B37', 'Maker'),
    @('C37', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'),
    @('D37', '1.168.35'),
    @('E37', '  +0.40%  '),
    @('E38', '  +2.45%  '),
    @('E39', '  -0.01%  '),
    @('D40', '0.505'),
    @('E40', '  +1.56%  '),
    @('E41', '  -0.23%  '),
    @('D42', '0.793'),
    @('E42', '  +1.19%  '),
    @('D43', '5.38'),
    @('E43', '  -0.15%  '),
    @('D44', '1.770.61'),
    @('E44', '  +1.48%  '),
    @('D45', '92.47'),
    @('D46', '1.56'),
    @('E46', '  +1.11%  '),
    @('D47', '54.68'),
    @('E47', '  +1.28%  '),
    @('E48', '  +0.85%  '),
    @('B49', 'EnergySwap'),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('D49', '7.59'),
    @('E49', '  +4.39%  '),
    @('B50', 'Mantle'),
    @('C50', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @('D50', '0.409'),
    @('E50', '  +0.36%  '),
    @('E51', '  -0.01%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $col = $cellRef.Substring(0, 1)
    $rng = $ws.Range($cellRef)
    if ($col -eq "D" -or $col -eq "E") {
        $rng.NumberFormat = "@"
        $rng.Value = $newVal
        $rng.Style = "Normal"
    } else {
        $rng.Value = $newVal
    }
}